$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Serie" dates for 04-08-2021 and 05-08-2021 must land in the sheet as
# plain text (shared strings), matching how every other date in column A is
# stored. Assigning a date-look-alike string straight to Range.Value makes
# Excel "smart" parse it into a date serial (it would pick up a number
# format / style that the source workbook never had). To avoid that, write
# the literal text with a leading apostrophe into a scratch cell far outside
# the used range (forces text, via quote-prefix), copy it, and paste
# *values only* into the target cell - a values-only paste does not carry
# the quote-prefix formatting over, so the destination cell ends up as a
# plain shared-string cell with no style override, exactly like the rest of
# column A.
$ws.Range("Z1").Value = "'04-08-2021"
$ws.Range("Z1").Copy()
$ws.Range("A56").PasteSpecial(-4163)

$ws.Range("Z1").Value = "'05-08-2021"
$ws.Range("Z1").Copy()
$ws.Range("A57").PasteSpecial(-4163)

$ws.Range("Z1").Clear()

$ws.Range("B56").Value = 200000
$ws.Range("C56").Value = 431000
$ws.Range("D56").Value = 200000
$ws.Range("E56").Value = 189000
$ws.Range("F56").Value = 11000
$ws.Range("G56").Value = 1

$ws.Range("B57").Value = 200000
$ws.Range("C57").Value = 321000
$ws.Range("D57").Value = 100000
$ws.Range("E57").Value = 89000
$ws.Range("F57").Value = 11000
$ws.Range("G57").Value = 1.03
